$d = $word.ActiveDocument
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# The document currently has a single paragraph that only contains the
# (hidden) "_GoBack" bookmark. That paragraph becomes the third paragraph
# ("Presentado por: ...") of the final document; we build the new title
# paragraph and the blank bold paragraph before it, and the four blank
# paragraphs + trailing paragraph after it.

$bm = $d.Paragraphs(1)

# --- Insert the two new paragraphs that precede the bookmark paragraph ---
# 1) "Actividad ejercicio prueba, preparatorio parcial " (bold)
# 2) an empty paragraph whose paragraph mark is bold
$beforeXml = "<w:p $ns>" +
  "<w:pPr><w:rPr><w:b/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:b/></w:rPr><w:t xml:space=`"preserve`">Actividad ejercicio prueba, preparatorio parcial </w:t></w:r>" +
  "</w:p>" +
  "<w:p $ns><w:pPr><w:rPr><w:b/></w:rPr></w:pPr></w:p>"

$insertBefore = $d.Range($bm.Range.Start, $bm.Range.Start)
$insertBefore.InsertXML($beforeXml)

# --- Turn the (still untouched) bookmark paragraph into the "Presentado
# por" paragraph, keeping the bookmark physically between the two runs ---
$bm = $d.Paragraphs(3)
$full = $d.Range($bm.Range.Start, $bm.Range.End)
$afterXml = "<w:p $ns><w:r><w:t>or: Mauricio Duque   cc 1036403902</w:t></w:r></w:p>"
$full.InsertXML($afterXml)

$bm = $d.Paragraphs(3)
$startPoint = $d.Range($bm.Range.Start, $bm.Range.Start)
$startPoint.InsertBefore("Presentado p")

$bm = $d.Paragraphs(3)
$bm.Range.Bold = 1

# --- Append the four blank paragraphs and the final paragraph (two
# spaces, not bold) after the "Presentado por" paragraph ---
$bm = $d.Paragraphs(3)
$tail = $bm.Range.InsertParagraphAfter()
$next = $d.Paragraphs(4)
$tailTarget = $d.Range($next.Range.Start, $next.Range.End)
$tailXml = "<w:p $ns/><w:p $ns/><w:p $ns/><w:p $ns/>" +
  "<w:p $ns><w:r><w:t xml:space=`"preserve`">  </w:t></w:r></w:p>"
$tailTarget.InsertXML($tailXml)
